# PlanilhaHorasEstagio_Alessandro.xlsx — add start/end time columns, a
# header row, and two new logbook rows (connecting to the DB + creating
# the entity), per commit "Conexão com banco e criação da entidade funci".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Split the single "time range" column into two columns (start/end).
#    Insert a fresh column before B; it pushes the old time-text column
#    (old B) to C, and the old merged description block (old C:J) to D:K.
# ------------------------------------------------------------------
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").ColumnWidth = 20

# ------------------------------------------------------------------
# 2) Header row (row 2 was blank before — just populate it).
# ------------------------------------------------------------------
$ws.Range("A2").Value = "DATA"
$ws.Range("B2").Value = "HORA INICIO"
$ws.Range("C2").Value = "HORA TÉRMINO"
$ws.Range("D2").Value = "DESCRIÇÃO ATIVIDADES"
$ws.Range("D2:K2").Merge()
$ws.Range("D2:K2").Style = $ws.Range("D3").Style

# ------------------------------------------------------------------
# 3) Replace the old "22:30 - 00:30" style text in column C (now shifted
#    there) with two real time-of-day values in B/C, time-formatted.
#    Column C keeps the old merge-cell style (s=2) carried over from the
#    description block start, so explicitly restamp B/C with the time
#    style and clear any stray leftovers.
# ------------------------------------------------------------------
$timeFmt = "h:mm:ss;@"

$ws.Range("B3").Value = 0.9375
$ws.Range("C3").Value = 0.020833333333333332

$ws.Range("B4").Value = 0.9375
$ws.Range("C4").Value = 0.020833333333333332

$ws.Range("B5").Value = 0.77083333333333337
$ws.Range("C5").Value = 0.9375

$ws.Range("B6").Value = 0.3125
$ws.Range("C6").Value = 0.52083333333333337

$ws.Range("B7").Value = 0.9375
$ws.Range("C7").Value = 0.10416666666666667

$ws.Range("B3:C7").NumberFormat = $timeFmt

# ------------------------------------------------------------------
# 4) Two brand-new log rows: connecting to the DB with JPA, and creating
#    new entities, plus a still-open in-progress row (date/time only).
# ------------------------------------------------------------------
$ws.Range("A8").Value = 44723
$ws.Range("B8").Value = 0.3125
$ws.Range("C8").Value = 0.58333333333333337
$ws.Range("D8").Value = "Conectando com o banco de dados e criando relação de persistência com JPA"
$ws.Range("D8:K8").Merge()
$ws.Range("D8:K8").Style = $ws.Range("D3").Style

$ws.Range("A9").Value = 44724
$ws.Range("B9").Value = 0.33333333333333331
$ws.Range("C9").Value = 0.041666666666666664
$ws.Range("D9").Value = "Criando novas entidades como logica de negócio"
$ws.Range("D9:K9").Merge()
$ws.Range("D9:K9").Style = $ws.Range("D3").Style

$ws.Range("A10").Value = 44731
$ws.Range("B10").Value = 0.91666666666666663
$ws.Range("C10").Value = 0.059027777777777783

$ws.Range("A3:A10").NumberFormat = "m/d/yyyy"
$ws.Range("B8:C10").NumberFormat = $timeFmt

# ------------------------------------------------------------------
# 5) View niceties matching the author's last on-screen state.
# ------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("C10").Select()
